{"js": "// Replace the text of each cell in the first (and only) table with its\n// new value, preserving existing run/paragraph formatting (font,\n// size, justification, etc.) by rewriting the paragraph's range text\n// instead of clearing+re-inserting the whole cell body.\n//\n// The new values below are the row-major (top-left -> bottom-right)\n// replacement sequence for the 20x5 table of arithmetic problems, built\n// from the canonical OOXML diff.\nconst newValues = [\n  \"7+19=26\", \"47-39=8\", \"57-29=28\", \"96-8=88\", \"62-14=48\",\n  \"52-5=47\", \"39+24=63\", \"9+38=47\", \"43-8=35\", \"53+29=82\",\n  \"41-2=39\", \"32-16=16\", \"62-18=44\", \"16+29=45\", \"6+69=75\",\n  \"65+7=72\", \"75-58=17\", \"12+79=91\", \"26+17=43\", \"9+27=36\",\n  \"54-7=47\", \"28+35=63\", \"74-25=49\", \"28-19=9\", \"44+48=92\",\n  \"13+68=81\", \"55+27=82\", \"66+29=95\", \"73-64=9\", \"67-59=8\",\n  \"29+23=52\", \"38+53=91\", \"15+27=42\", \"18+49=67\", \"60-53=7\",\n  \"59+5=64\", \"92-19=73\", \"72-24=48\", \"92-27=65\", \"7+87=94\",\n  \"45-9=36\", \"92-5=87\", \"90-76=14\", \"78+15=93\", \"19+39=58\",\n  \"57+7=64\", \"38+13=51\", \"19+69=88\", \"58-49=9\", \"37+37=74\",\n  \"68+17=85\", \"81-55=26\", \"86-49=37\", \"72-58=14\", \"13+78=91\",\n  \"68+13=81\", \"87+6=93\", \"80-15=65\", \"5+7=12\", \"18+28=46\",\n  \"7+75=82\", \"7+36=43\", \"88+7=95\", \"54-27=27\", \"38+25=63\",\n  \"54-8=46\", \"41-25=16\", \"36+57=93\", \"55-28=27\", \"25+9=34\",\n  \"93-45=48\", \"73-9=64\", \"41-35=6\", \"22-9=13\", \"67+5=72\",\n  \"8+56=64\", \"92-28=64\", \"5+9=14\", \"96-67=29\", \"73+18=91\",\n  \"82-43=39\", \"60-12=48\", \"6+79=85\", \"46+28=74\", \"70-6=64\",\n  \"24+28=52\", \"7+89=96\", \"7+9=16\", \"85-38=47\", \"25-9=16\",\n  \"36+59=95\", \"50-45=5\", \"47+18=65\", \"70-58=12\", \"46+49=95\",\n  \"88-19=69\", \"67-8=59\", \"38+3=41\", \"46+18=64\", \"90-33=57\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values[0].length;\n\nif (rowCount * colCount !== newValues.length) {\n  throw new Error(\n    `Expected ${newValues.length} cells, found ${rowCount}x${colCount}=${rowCount * colCount}`\n  );\n}\n\n// Collect every cell's first paragraph so we can batch the loads.\nconst paragraphs = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    paragraphs.push(cell.body.paragraphs);\n  }\n}\nawait context.sync();\n\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const para = paragraphs[i].items[0];\n    para.getRange().insertText(newValues[i], Word.InsertLocation.replace);\n    i++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the text of each cell in the first (and only) table with its\n# new value, preserving existing run/paragraph formatting (font, size,\n# justification, etc.) by setting Cell.Range.Text directly instead of\n# deleting/recreating cell content.\n#\n# The values below are the row-major (top-left -> bottom-right)\n# replacement sequence for the 20x5 table of arithmetic problems, built\n# from the canonical OOXML diff.\n$newValues = @(\n  \"7+19=26\",\"47-39=8\",\"57-29=28\",\"96-8=88\",\"62-14=48\",\n  \"52-5=47\",\"39+24=63\",\"9+38=47\",\"43-8=35\",\"53+29=82\",\n  \"41-2=39\",\"32-16=16\",\"62-18=44\",\"16+29=45\",\"6+69=75\",\n  \"65+7=72\",\"75-58=17\",\"12+79=91\",\"26+17=43\",\"9+27=36\",\n  \"54-7=47\",\"28+35=63\",\"74-25=49\",\"28-19=9\",\"44+48=92\",\n  \"13+68=81\",\"55+27=82\",\"66+29=95\",\"73-64=9\",\"67-59=8\",\n  \"29+23=52\",\"38+53=91\",\"15+27=42\",\"18+49=67\",\"60-53=7\",\n  \"59+5=64\",\"92-19=73\",\"72-24=48\",\"92-27=65\",\"7+87=94\",\n  \"45-9=36\",\"92-5=87\",\"90-76=14\",\"78+15=93\",\"19+39=58\",\n  \"57+7=64\",\"38+13=51\",\"19+69=88\",\"58-49=9\",\"37+37=74\",\n  \"68+17=85\",\"81-55=26\",\"86-49=37\",\"72-58=14\",\"13+78=91\",\n  \"68+13=81\",\"87+6=93\",\"80-15=65\",\"5+7=12\",\"18+28=46\",\n  \"7+75=82\",\"7+36=43\",\"88+7=95\",\"54-27=27\",\"38+25=63\",\n  \"54-8=46\",\"41-25=16\",\"36+57=93\",\"55-28=27\",\"25+9=34\",\n  \"93-45=48\",\"73-9=64\",\"41-35=6\",\"22-9=13\",\"67+5=72\",\n  \"8+56=64\",\"92-28=64\",\"5+9=14\",\"96-67=29\",\"73+18=91\",\n  \"82-43=39\",\"60-12=48\",\"6+79=85\",\"46+28=74\",\"70-6=64\",\n  \"24+28=52\",\"7+89=96\",\"7+9=16\",\"85-38=47\",\"25-9=16\",\n  \"36+59=95\",\"50-45=5\",\"47+18=65\",\"70-58=12\",\"46+49=95\",\n  \"88-19=69\",\"67-8=59\",\"38+3=41\",\"46+18=64\",\"90-33=57\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nif (($rowCount * $colCount) -ne $newValues.Length) {\n  throw \"Expected $($newValues.Length) cells, found $rowCount x $colCount = $($rowCount * $colCount)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$i]\n    $i++\n  }\n}\n"}
